# Working on form data verification
# - Update the submitted e-mail value in AK2 so it includes the full
#   address (matches what is already used as the hyperlink's mailto target).
# - Re-apply the hyperlink so its "display" text is no longer redundantly
#   stored (Excel omits `display` when it matches the cell text).
# - Move the active-cell selection on the sheet to AE12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCell = $ws.Range("AK2")
$emailAddress = "sindhube19.data@gmail.com"
$mailTarget = "mailto:" + $emailAddress

# Drop the existing hyperlink (and its now-redundant display text), then
# update the cell's text and re-insert the hyperlink pointing at the same
# mailto address so the display text is implicitly the cell value again.
$targetCell.Hyperlinks.Delete()
$targetCell.Value = $emailAddress
$ws.Hyperlinks.Add($targetCell, $mailTarget)

# Restore the built-in "Hyperlink" cell style (re-adding the hyperlink
# re-applies formatting through a fresh style record; putting the named
# style back keeps the cell on the original shared style).
$targetCell.Style = "Hyperlink"

# Update the current selection to AE12.
$ws.Range("AE12").Select()
